# Generate Report for Handoff
# Update localization-status workbook: mark b.md as "Ready for handoff"
# on all sheets, and record the new handoff file / datetime / error detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 14:45:48"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"

# "False" (and "True") are auto-detected as booleans by Value, so enter the
# text via a formula and convert the formula cell to a literal value in
# place, which keeps it typed as a shared string instead of xlsx boolean.
$wsZhCn.Range("F3").Formula = '="False"'
$wsZhCn.Range("F3").Copy()
$wsZhCn.Range("F3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 14:45:44"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a92e8a194ea4d986cba6b9ead9572cecd26361f8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2201cce3dfb6049a23153fcaf2f30bb4cef654a7/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"

$wsDeDe.Range("F3").Formula = '="False"'
$wsDeDe.Range("F3").Copy()
$wsDeDe.Range("F3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-21 14:45:48"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a92e8a194ea4d986cba6b9ead9572cecd26361f8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2201cce3dfb6049a23153fcaf2f30bb4cef654a7/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 40
